# 自动更新Excel文件 - daily refresh of remaining-days (剩余) counters.
#
# Column layout (row 1 = header):
#   A 行号    B 店铺名称   C 地址   D 总天 (total days)
#   E 剩余 (remaining days)   F 开始时间 (start date, stored as an
#       integer in YYYYMMDD form, NOT a real Excel date serial)
#   G/H/I 备注1..3
#
# Business rule recovered from the data: remaining = D - (today - F),
# where "today" is the reference date implied by the sheet (it advances
# by one calendar day on every run: 2025-11-20 -> 2025-11-21 here).
# So on every refresh each row's E simply decrements by 1. The one
# exception is a row whose countdown would hit 0: instead of rolling to
# 0, the row's "start date" is reset to the new "today" and its
# remaining count is reset back to its full total (D) - i.e. the cycle
# restarts.
#
# Rows with an unparseable start date (data-entry typo) are left
# untouched, exactly like the source diff leaves row 36 alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- pure-integer civil-calendar <-> day-number helpers -------------
# (Howard Hinnant's days_from_civil / civil_from_days, integer only —
# avoids any reliance on .NET DateTime/TimeSpan arithmetic.)

function Get-DaysFromCivil($y, $m, $d) {
    $y2 = $y
    if ($m -le 2) { $y2 = $y2 - 1 }
    if ($y2 -ge 0) { $eraBase = $y2 } else { $eraBase = $y2 - 399 }
    $era = [math]::Floor($eraBase / 400)
    $yoe = $y2 - ($era * 400)
    $mp = (($m + 9) % 12)
    $doy = [math]::Floor((153 * $mp + 2) / 5) + $d - 1
    $doe = ($yoe * 365) + [math]::Floor($yoe / 4) - [math]::Floor($yoe / 100) + $doy
    return ($era * 146097) + $doe - 719468
}

function Get-CivilFromDays($z) {
    $z2 = $z + 719468
    if ($z2 -ge 0) { $eraBase = $z2 } else { $eraBase = $z2 - 146096 }
    $era = [math]::Floor($eraBase / 146097)
    $doe = $z2 - ($era * 146097)
    $yoe = [math]::Floor(($doe - [math]::Floor($doe / 1460) + [math]::Floor($doe / 36524) - [math]::Floor($doe / 146096)) / 365)
    $y = $yoe + ($era * 400)
    $doy = $doe - ((365 * $yoe) + [math]::Floor($yoe / 4) - [math]::Floor($yoe / 100))
    $mp = [math]::Floor((5 * $doy + 2) / 153)
    $d = $doy - [math]::Floor((153 * $mp + 2) / 5) + 1
    if ($mp -lt 10) { $m = $mp + 3 } else { $m = $mp - 9 }
    if ($m -le 2) { $y = $y + 1 }
    return @($y, $m, $d)
}

function ConvertTo-YMD($ymdInt) {
    # ymdInt like 20251117 -> (2025, 11, 17); $null if not a clean
    # 8-digit yyyymmdd number.
    if ($ymdInt -eq $null) { return $null }
    $n = [math]::Round($ymdInt)
    if ($n -lt 10000101 -or $n -gt 99991231) { return $null }
    $y = [math]::Floor($n / 10000)
    $m = [math]::Floor(($n % 10000) / 100)
    $d = $n % 100
    if ($m -lt 1 -or $m -gt 12 -or $d -lt 1 -or $d -gt 31) { return $null }
    return @($y, $m, $d)
}

# ---- locate the reference "today" from the sheet's current state ----
# Every (D,E,F) triple should currently imply the same "today"; use the
# first row that parses cleanly to establish it, then advance by 1 day.

$dim = $ws.UsedRange
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

$oldToday = $null
for ($r = 2; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $fVal = $ws.Cells.Item($r, 6).Value2
    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) { continue }
    $ymd = ConvertTo-YMD $fVal
    if ($ymd -eq $null) { continue }
    $startDays = Get-DaysFromCivil $ymd[0] $ymd[1] $ymd[2]
    $elapsed = $dVal - $eVal
    $oldToday = $startDays + $elapsed
    break
}

if ($oldToday -eq $null) {
    # Fallback: nothing usable on the sheet, nothing to refresh.
    return
}

$newToday = $oldToday + 1
$newTodayYmd = Get-CivilFromDays $newToday
$newTodayInt = ($newTodayYmd[0] * 10000) + ($newTodayYmd[1] * 100) + $newTodayYmd[2]

# ---- refresh every data row ------------------------------------------

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2
    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) { continue }

    $ymd = ConvertTo-YMD $fVal
    if ($ymd -eq $null) { continue }

    $newRemaining = $eVal - 1
    if ($newRemaining -le 0) {
        # countdown finished -> restart the cycle from "today"
        $fCell.Value = $newTodayInt
        $eCell.Value = $dVal
    } else {
        $eCell.Value = $newRemaining
    }
}
